$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "20.743.54"
$ws.Range("E2").Value = "  -5.88%  "
$ws.Range("D3").Value = "1.472.07"
$ws.Range("E3").Value = "  -5.42%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("E4").Value = "  +0.96%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.008"
$ws.Range("E5").Value = "  +0.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "279.63"
$ws.Range("E6").Value = "  -3.82%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3758"
$ws.Range("E7").Value = "  -5.40%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3080"
$ws.Range("E8").Value = "  -4.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "41.99"
$ws.Range("E9").Value = "  -5.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06720"
$ws.Range("E10").Value = "  -7.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.013"
$ws.Range("E11").Value = "  -6.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.010"
$ws.Range("E12").Value = "  +1.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.470"
$ws.Range("E13").Value = "  -4.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.60"
$ws.Range("E14").Value = "  -6.44%  "
$ws.Range("D15").Value = "1.488.89"
$ws.Range("E15").Value = "  -4.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.238"
$ws.Range("E16").Value = "  -6.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001031"
$ws.Range("E17").Value = "  -8.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06480"
$ws.Range("E18").Value = "  -1.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "79.67"
$ws.Range("E19").Value = "  -4.81%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.008"
$ws.Range("E20").Value = "  +0.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.848"
$ws.Range("E21").Value = "  -6.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.78"
$ws.Range("E22").Value = "  -5.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.83"
$ws.Range("E23").Value = "  -4.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.327"
$ws.Range("E24").Value = "  -1.56%  "
$ws.Range("D25").Value = "20.748.89"
$ws.Range("E25").Value = "  -5.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.273"
$ws.Range("E26").Value = "  -6.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "145.10"
$ws.Range("E27").Value = "  -2.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.58"
$ws.Range("E28").Value = "  -5.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.775"
$ws.Range("E29").Value = "  -1.96%  "
$ws.Range("D30").Value = "1.655.97"
$ws.Range("E30").Value = "  -4.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "111.55"
$ws.Range("E31").Value = "  -6.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.696"
$ws.Range("E32").Value = "  -2.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9200"
$ws.Range("E33").Value = "  -7.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07848"
$ws.Range("E34").Value = "  -5.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.371"
$ws.Range("E35").Value = "  -8.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.455"
$ws.Range("E36").Value = "  -9.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.12"
$ws.Range("E37").Value = "  +3.30%  "
$ws.Range("B38").Value = "Frax"
$ws.Range("C38").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.007"
$ws.Range("E38").Value = "  +0.83%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.849"
$ws.Range("E39").Value = "  -5.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05722"
$ws.Range("E40").Value = "  -4.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1939"
$ws.Range("E41").Value = "  -5.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.02067"
$ws.Range("E42").Value = "  -8.73%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.129"
$ws.Range("E43").Value = "  -6.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5454"
$ws.Range("E44").Value = "  -6.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.622"
$ws.Range("E45").Value = "  -3.31%  "
$ws.Range("E46").Value = "  -4.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5262"
$ws.Range("E47").Value = "  -5.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.807"
$ws.Range("E48").Value = "  -5.01%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "110.89"
$ws.Range("E49").Value = "  -6.33%  "
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.090"
$ws.Range("E50").Value = "  -4.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06398"
$ws.Range("E51").Value = "  -6.35%  "
